$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 32; $r -le 35; $r++) {
    $ws.Cells.Item($r, 9).Value = "Approved"
    $ws.Cells.Item($r, 10).ClearContents()
}

$ws.Range("I31").Select()
